$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 197.5
$ws.Range("I41").Value = 210
$ws.Range("J41").Value = 160
$ws.Range("K41").Value = 210
$ws.Range("L41").Value = 160
$ws.Range("M41").Value = 230
$ws.Range("N41").Value = -1040

$ws.Range("H53").Value = 257.625
$ws.Range("I53").Value = 290.25
$ws.Range("J53").Value = 225
$ws.Range("K53").Value = 290.25
$ws.Range("L53").Value = 225
$ws.Range("M53").Value = 346.75
$ws.Range("N53").Value = -1499

$ws.Range("H113").Value = 2875
$ws.Range("I113").Value = 2700
$ws.Range("J113").Value = 3050
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 3050
$ws.Range("M113").Value = 554
$ws.Range("N113").Value = -9558

$ws.Range("H125").Value = 1749.5
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1749.5
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 15745.5
$ws.Range("N125").Value = -20665.5
$ws.Range("M125").Value = $null

$ws.Range("H137").Value = 2355.65
$ws.Range("I137").Value = 1966.3438
$ws.Range("J137").Value = 3912.875
$ws.Range("K137").Value = 5899.0314
$ws.Range("L137").Value = 11738.625
$ws.Range("M137").Value = -3349.0314
$ws.Range("N137").Value = -16838.625

$ws.Range("H138").Value = 2167.2834
$ws.Range("J138").Value = 2632.9375
$ws.Range("L138").Value = 7898.8125
$ws.Range("N138").Value = -18178.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 669523.9399999999
$ws.Range("I32").Value = 765105.9399999999
$ws.Range("K32").Value = 765105.9399999999
$ws.Range("M32").Value = -764818.9399999999

$ws.Range("H74").Value = 806.7
$ws.Range("I74").Value = 682.63336
$ws.Range("K74").Value = 682.63336
$ws.Range("M74").Value = 191.36664

$ws.Range("H77").Value = 806.7
$ws.Range("I77").Value = 682.63336
$ws.Range("K77").Value = 3413.1668
$ws.Range("M77").Value = 954.8332

$ws.Range("H97").Value = 903.3333
$ws.Range("I97").Value = 903.3333
$ws.Range("K97").Value = 903.3333
$ws.Range("M97").Value = -407.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1220.1578
$ws.Range("I99").Value = 1204.6111
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1204.6111
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 293.3888999999999
$ws.Range("N99").Value = -4496

$ws.Range("H107").Value = 168916.67
$ws.Range("I107").Value = 251875
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 251875
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -249955
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 1861.1277
$ws.Range("I134").Value = 1476.7715
$ws.Range("J134").Value = 2982.1667
$ws.Range("K134").Value = 4430.3145
$ws.Range("L134").Value = 8946.500100000001
$ws.Range("M134").Value = -1895.3145
$ws.Range("N134").Value = -14016.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4210.9385
$ws.Range("I31").Value = 1052.2046
$ws.Range("K31").Value = 1052.2046
$ws.Range("M31").Value = -757.2046

$ws.Range("H34").Value = 4210.9385
$ws.Range("I34").Value = 1052.2046
$ws.Range("K34").Value = 1052.2046
$ws.Range("M34").Value = -850.2046

$ws.Range("H58").Value = 1333.0312
$ws.Range("I58").Value = 1159.7727
$ws.Range("K58").Value = 1159.7727
$ws.Range("M58").Value = -956.7727

$ws.Range("H107").Value = 6251150
$ws.Range("I107").Value = 8929500
$ws.Range("J107").Value = 1666.6666
$ws.Range("K107").Value = 8929500
$ws.Range("L107").Value = 1666.6666
$ws.Range("M107").Value = -8927580
$ws.Range("N107").Value = -5506.6666

$ws.Range("H136").Value = 1333.0312
$ws.Range("I136").Value = 1159.7727
$ws.Range("K136").Value = 3479.3181
$ws.Range("M136").Value = -929.3181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 904.6842
$ws.Range("J5").Value = 1972.2
$ws.Range("L5").Value = 5916.6
$ws.Range("N5").Value = -6140.6

$ws.Range("H132").Value = 2158.8333
$ws.Range("I132").Value = 1620.75
$ws.Range("J132").Value = 2427.875
$ws.Range("K132").Value = 14586.75
$ws.Range("L132").Value = 21850.875
$ws.Range("M132").Value = -12056.75
$ws.Range("N132").Value = -26910.875

$ws.Range("H135").Value = 904.6842
$ws.Range("J135").Value = 1972.2
$ws.Range("L135").Value = 17749.8
$ws.Range("N135").Value = -22819.8

$ws.Range("H137").Value = 8343456
$ws.Range("I137").Value = 55593892
$ws.Range("J137").Value = 5143.7646
$ws.Range("K137").Value = 166781676
$ws.Range("L137").Value = 15431.2938
$ws.Range("M137").Value = -166776576
$ws.Range("N137").Value = -25631.2938

$ws.Range("H138").Value = 7166
$ws.Range("J138").Value = 11000
$ws.Range("L138").Value = 33000
$ws.Range("N138").Value = -43280

$ws.Range("H140").Value = 1202.7858
$ws.Range("I140").Value = 1202.7858
$ws.Range("K140").Value = 3608.3574
$ws.Range("M140").Value = 1571.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 31000
$ws.Range("J62").Value = 31000
$ws.Range("L62").Value = 31000
$ws.Range("N62").Value = -32372

$ws.Range("H65").Value = 31000
$ws.Range("J65").Value = 31000
$ws.Range("L65").Value = 93000
$ws.Range("N65").Value = -99864

$ws.Range("H97").Value = 1673.8889
$ws.Range("I97").Value = 1233
$ws.Range("J97").Value = 2225
$ws.Range("K97").Value = 1233
$ws.Range("L97").Value = 2225
$ws.Range("M97").Value = -737
$ws.Range("N97").Value = -3217

$ws.Range("H122").Value = 20001458
$ws.Range("I122").Value = 1539.1818
$ws.Range("J122").Value = 166667540
$ws.Range("K122").Value = 4617.5454
$ws.Range("L122").Value = 500002620
$ws.Range("M122").Value = -2167.5454
$ws.Range("N122").Value = -500007520

$ws.Range("H127").Value = 55155.168
$ws.Range("J127").Value = 55155.168
$ws.Range("L127").Value = 55155.168
$ws.Range("N127").Value = -65075.168

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H132").Value = 1896.138
$ws.Range("I132").Value = 1582.5
$ws.Range("J132").Value = 3401.6
$ws.Range("K132").Value = 4747.5
$ws.Range("L132").Value = 10204.8
$ws.Range("M132").Value = -2217.5
$ws.Range("N132").Value = -15264.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5761.346
$ws.Range("I22").Value = 1363.6364
$ws.Range("J22").Value = 8986.333000000001
$ws.Range("K22").Value = 1363.6364
$ws.Range("L22").Value = 8986.333000000001
$ws.Range("M22").Value = -1068.6364
$ws.Range("N22").Value = -9576.333000000001

$ws.Range("H27").Value = 5761.346
$ws.Range("I27").Value = 1363.6364
$ws.Range("J27").Value = 8986.333000000001
$ws.Range("K27").Value = 1363.6364
$ws.Range("L27").Value = 8986.333000000001
$ws.Range("M27").Value = -1256.6364
$ws.Range("N27").Value = -9200.333000000001

$ws.Range("H40").Value = 169317.33
$ws.Range("I40").Value = 252226
$ws.Range("K40").Value = 252226
$ws.Range("M40").Value = -252090

$ws.Range("H61").Value = 4300
$ws.Range("I61").Value = 1600
$ws.Range("J61").Value = 4840
$ws.Range("K61").Value = 1600
$ws.Range("L61").Value = 4840
$ws.Range("M61").Value = -1398
$ws.Range("N61").Value = -5244

$ws.Range("H113").Value = 4300
$ws.Range("I113").Value = 1600
$ws.Range("J113").Value = 4840
$ws.Range("K113").Value = 1600
$ws.Range("L113").Value = 4840
$ws.Range("M113").Value = 570
$ws.Range("N113").Value = -9180

$ws.Range("H122").Value = 3522.7368
$ws.Range("I122").Value = 3330.5557
$ws.Range("J122").Value = 3695.7
$ws.Range("K122").Value = 9991.667099999999
$ws.Range("L122").Value = 11087.1
$ws.Range("M122").Value = -7541.667099999999
$ws.Range("N122").Value = -15987.1

$ws.Range("H136").Value = 8773644
$ws.Range("I136").Value = 1815.5454
$ws.Range("J136").Value = 20834908
$ws.Range("K136").Value = 5446.6362
$ws.Range("L136").Value = 62504724
$ws.Range("M136").Value = -2896.6362
$ws.Range("N136").Value = -62509824
